# "new power plants and preparation script"
#   - rows 2-29 (existing plants) get revised Technology/Age/Capacity/
#     Efficiency/availability_weighted_mean values
#   - rows 30-51 are brand-new plant entries appended below the table
#     (dimension grows from A1:J29 to A1:J51)

function Write-PlantRows($Sheet, $StartRow, $Rows) {
    # $Rows is an array of 9-element arrays: Technology, Age, Capacity,
    # Efficiency, availability_weighted_mean, Location, Owner, cash,
    # DischarginEfficiency (columns B..J)
    $arr = New-Object 'object[,]' $Rows.Count, 9
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        for ($j = 0; $j -lt 9; $j++) {
            $arr[$i, $j] = $Rows[$i][$j]
        }
    }
    $endRow = $StartRow + $Rows.Count - 1
    $Sheet.Range("B$StartRow`:J$endRow").Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1) revised data for the existing rows 2-29 (row id in col A is unchanged) --
$existingRows = @(
    ,@("Coal PSC", 4, 790, 0.46, 0.9, "NL", "Producer1", 0, 0)
    ,@("Coal PSC", 5, 2640, 0.46, 0.9, "NL", "Producer1", 0, 0)
    ,@("Coal PSC", 24, 630, 0.41, 0.9, "NL", "Producer1", 0, 0)
    ,@("Coal PSC", 25, 631, 0.42, 0.8000000000000002, "NL", "Producer1", 0, 0)
    ,@("CCGT", 1, 341, 0.52, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 2, 800, 0.5600000000000001, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 3, 874, 0.58, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 4, 420, 0.59, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 6, 872, 0.5849885321100917, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 7, 420, 0.59, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 8, 426, 0.58, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 9, 1295, 0.587984555984556, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 10, 870, 0.59, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 19, 245, 0.51, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 22, 682, 0.52, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT", 23, 341, 0.52, 0.9, "NL", "Producer1", 0, 0)
    ,@("Hydropower_reservoir_medium", 14, 40.207, 1, 1, "NL", "Producer1", 0, 0)
    ,@("Nuclear", 13, 484, 0.34, 0.94999999975, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 0, 2400.537, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 1, 1527.304, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 2, 588.5690000000001, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 3, 633.914, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 4, 498.2, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 5, 492.642, 1, 1, "NL", "Producer1", 0, 0)
    ,@("PV_utility_systems", 6, 648.744, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_offshore", 2, 203.696, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_offshore", 3, 758.994, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 0, 341.774, 1, 1, "NL", "Producer1", 0, 0)
)
Write-PlantRows $ws 2 $existingRows

# -- 2) brand-new plants appended as rows 30-51 --
$newRowIds = @(28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49)
$newRows = @(
    ,@("WTG_onshore", 1, 425.903, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 2, 551.907, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 3, 659.4730000000001, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 11, 638.6619999999999, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 12, 317.364, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 13, 285.598, 1, 1, "NL", "Producer1", 0, 0)
    ,@("WTG_onshore", 14, 1016.313, 1, 1, "NL", "Producer1", 0, 0)
    ,@("OCGT", 33, 643, 0.387651632970451, 0.9, "NL", "Producer1", 0, 0)
    ,@("OCGT", 45, 460, 0.363, 0.9, "NL", "Producer1", 0, 0)
    ,@("Biomass_CHP_wood_pellets_PH", 14, 1.767999999999997, 0.3, 0.742009132420091, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 6, 444.2, 0.5844079243583972, 0.8917154434939217, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 11, 862.64, 0.4099999999999999, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 12, 721, 0.4099999999999999, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 13, 641, 0.41, 0.9, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 14, 739.293, 0.2941889890070649, 0.8515333832623808, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 19, 541.817, 0.4088526199805469, 0.804654952562272, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 20, 557.934, 0.4158746375019268, 0.9035726544813641, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 22, 1040.24, 0.3575534876566946, 0.90748031496063, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 23, 1101.962, 0.4374676880872479, 0.8645218458763383, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 28, 655.18101, 0.352119683841264, 0.8390035408294167, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 29, 385.096, 0.3745767289195421, 0.8420443214159588, "NL", "Producer1", 0, 0)
    ,@("CCGT_CHP_backpressure_DH", 32, 811.6, 0.3692304313085264, 0.8862000985707243, "NL", "Producer1", 0, 0)
)

# stamp column A: copy the existing header-row style (bold/bordered) down
# into A30:A51, then fill in the row ids
$ws.Range("A2").Copy($ws.Range("A30:A51"))
$idArr = New-Object 'object[,]' $newRowIds.Count, 1
for ($i = 0; $i -lt $newRowIds.Count; $i++) { $idArr[$i, 0] = $newRowIds[$i] }
$ws.Range("A30:A51").Value = $idArr
Write-PlantRows $ws 30 $newRows

